$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (FAPs -> Inflammatory-Mac) ---
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("I2").Value = 0.4877525841056716
$ws.Range("J2").Value = 0.588184597482006
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.375733333333333
$ws.Range("N2").Value = 4.1272
$ws.Range("O2").Value = 0.457732955319909
$ws.Range("P2").Value = 0.457732955319909
$ws.Range("Q2").Value = 0.5192123072888889
$ws.Range("R2").Value = 4.6729107656
$ws.Range("S2").Value = 0.2232604317876115
$ws.Range("T2").Value = 0.2692314740790897

# --- Row 3 (FAPs -> Resolving-Mac) ---
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("I3").Value = 0.4877525841056716
$ws.Range("J3").Value = 0.588184597482006
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.629803666666667
$ws.Range("N3").Value = 4.889411000000001
$ws.Range("O3").Value = 0.542267044680091
$ws.Range("P3").Value = 0.542267044680091
$ws.Range("Q3").Value = 0.6151003989614445
$ws.Range("R3").Value = 5.535903590653001
$ws.Range("S3").Value = 0.2644921523180601
$ws.Range("T3").Value = 0.3189531234029163

# --- Row 4 (MuSCs -> Inflammatory-Mac) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.396361
$ws.Range("H4").Value = 0.7927219999999999
$ws.Range("I4").Value = 0.5122474158943284
$ws.Range("J4").Value = 0.411815402517994
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.375733333333333
$ws.Range("N4").Value = 4.1272
$ws.Range("O4").Value = 0.457732955319909
$ws.Range("P4").Value = 0.457732955319909
$ws.Range("Q4").Value = 0.5452870397333334
$ws.Range("R4").Value = 3.2717222384
$ws.Range("S4").Value = 0.2344725235322975
$ws.Range("T4").Value = 0.1885014812408193

# --- Row 5 (MuSCs -> Resolving-Mac) ---
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.396361
$ws.Range("H5").Value = 0.7927219999999999
$ws.Range("I5").Value = 0.5122474158943284
$ws.Range("J5").Value = 0.411815402517994
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.629803666666667
$ws.Range("N5").Value = 4.889411000000001
$ws.Range("O5").Value = 0.542267044680091
$ws.Range("P5").Value = 0.542267044680091
$ws.Range("Q5").Value = 0.6459906111236667
$ws.Range("R5").Value = 3.875943666742
$ws.Range("S5").Value = 0.277774892362031
$ws.Range("T5").Value = 0.2233139212771747

# Remove the now-obsolete rows 6 and 7 (old MuSCs/Inflammatory-Mac, MuSCs/Resolving-Mac rows)
$ws.Range("A6:A7").EntireRow.Delete()
